$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1252
$ws.Range("F4").Value = 896
$ws.Range("F5").Value = 921
$ws.Range("F6").Value = 1665
$ws.Range("F7").Value = 362
$ws.Range("F8").Value = 1118
$ws.Range("F10").Value = 97
$ws.Range("F11").Value = 244
$ws.Range("F12").Value = 10
$ws.Range("F14").Value = 601
$ws.Range("F15").Value = 115
$ws.Range("F16").Value = 71
$ws.Range("F20").Value = 64
$ws.Range("F21").Value = 627
$ws.Range("F22").Value = 616
$ws.Range("F23").Value = 104
$ws.Range("F25").Value = 824
$ws.Range("F26").Value = 285
$ws.Range("F27").Value = 43
$ws.Range("F28").Value = 18
$ws.Range("F29").Value = 233
$ws.Range("F31").Value = 4
$ws.Range("F32").Value = 391

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 6
$ws.Range("F7").Value = 226
$ws.Range("F8").Value = 80

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 292

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 292
$ws.Range("F4").Value = 1252
$ws.Range("F5").Value = 896
$ws.Range("F6").Value = 921
$ws.Range("F7").Value = 1665
$ws.Range("F8").Value = 362
$ws.Range("F9").Value = 1118
$ws.Range("F12").Value = 97
$ws.Range("F13").Value = 244
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 601
$ws.Range("F17").Value = 115
$ws.Range("F18").Value = 71
$ws.Range("F23").Value = 6
$ws.Range("F26").Value = 226
$ws.Range("F27").Value = 226
$ws.Range("F28").Value = 64
$ws.Range("F29").Value = 627
$ws.Range("F30").Value = 616
$ws.Range("F31").Value = 104
$ws.Range("F33").Value = 824
$ws.Range("F34").Value = 285
$ws.Range("F35").Value = 80
$ws.Range("F36").Value = 43
$ws.Range("F37").Value = 18
$ws.Range("F38").Value = 233
$ws.Range("F43").Value = 4
$ws.Range("F45").Value = 391
